$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.748.86"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "1.864.12"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("D4").Value = "'1.034"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'323.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "'0.4425"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").Value = "'0.3804"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.07473"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").Value = "'0.8873"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").Value = "'21.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").Value = "1.874.62"
$ws.Range("E12").Value = "  -5.76%  "
$ws.Range("D13").Value = "'5.549"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "'0.07206"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "'84.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "'1.040"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "'0.000009147"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'1.032"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Value = "27.769.99"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").Value = "'5.320"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").Value = "'11.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "2.091.80"
$ws.Range("E24").Value = "  -3.95%  "
$ws.Range("E25").Value = "  +6.31%  "
$ws.Range("D26").Value = "'158.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").Value = "'18.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'5.382"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "'1.992"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("D30").Value = "'119.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.17%  "
$ws.Range("D31").Value = "'0.09068"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").Value = "'1.231"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("D33").Value = "'0.7798"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("D34").Value = "'3.035"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.05%  "
$ws.Range("D35").Value = "'4.603"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.98%  "
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").Value = "'1.144"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "'0.01992"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").Value = "'0.05363"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.36%  "
$ws.Range("D40").Value = "'2.884"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("D41").Value = "'0.5222"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("D42").Value = "'0.1697"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").Value = "'6.917"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.08%  "
$ws.Range("D44").Value = "'8.726"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("D45").Value = "'110.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.06696"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.40%  "
$ws.Range("D48").Value = "'1.037"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D49").Value = "'1.719"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("D50").Value = "'0.4736"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("E51").Value = "  +1.73%  "
